$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.479.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.102.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.095.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.72%  "
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.604.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.444.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.112"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.097.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "508.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.16%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("E30").Value = "  -10.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.43%  "
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "525.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.43%  "
$ws.Range("E38").Value = "  -3.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.079.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0795"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.121"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.254"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("E45").Value = "  +76.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₃0510"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.47%  "
